$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "normal effect" (2;21300101;1;100;1, group id 2 = buff/effect id 21300101)
# is removed from the gift-item list stored in B4 ("22031001" gift), matching
# issue #44 "finish the normal effect support".
$ws.Range("B4").Value = "1;22033001;1;100;5|1;22033002;1;100;5|1;22032007;1;100;1|1;22033013;1;100;5|1;22033014;1;100;3|1;22033015;1;100;3|2;21200101;1;100;1|2;21400101;1;100;1"

# Move the active selection to B4, matching the saved cursor position in the
# edited workbook.
$ws.Range("B4").Select() | Out-Null
